$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C3").Value = "Jashanpreet Sidhu"
$ws.Range("E7").Value = "None"
$ws.Range("F7").Value = "account_number = 2000`nclient_number = 2000`nbalance = 2000`ndate_created = (2015, 1, 1)`noverdraft_limit = -20`noverdraft_rate = 0.06`n"
$ws.Range("G7").Value = "attributes are set "
$ws.Range("E8").Value = "None"
$ws.Range("F8").Value = "account_number = 2000`nclient_number = 2000`nbalance = 2000`ndate_created = (2015, 1, 1)`noverdraft_limit = `"twenty`"`noverdraft_rate = 0.06"
$ws.Range("G8").Value = "overdraft_limit set to -100"
$ws.Range("E9").Value = "None"
$ws.Range("F9").Value = "account_number = 2000`nclient_number = 2000`nbalance = 2000`ndate_created = (2015, 1, 1)`noverdraft_limit = -20`noverdraft_rate = '6%'"
$ws.Range("G9").Value = "overdraft_rate set to 0.05"
$ws.Range("E10").Value = "None"
$ws.Range("F10").Value = "account_number = 2000`nclient_number = 2000`nbalance = 2000`ndate_created = '(1/1/15)'`noverdraft_limit = -20`noverdraft_rate = 0.06"
$ws.Range("G10").Value = "date_created set to today"
$ws.Range("E11").Value = "None"
$ws.Range("F11").Value = "account_number = 2000`nclient_number = 2000`nbalance = 2000`ndate_created = (2015, 1, 1)`noverdraft_limit = -20`noverdraft_rate = 0.06"
$ws.Range("G11").Value = "service_charge set to base service charge "
$ws.Range("E12").Value = "None"
$ws.Range("F12").Value = "account_number = 2000`nclient_number = 2000`nbalance = -2000`ndate_created = (2015, 1, 1)`noverdraft_limit = -20`noverdraft_rate = 0.06"
$ws.Range("G12").Value = "service_charge calculated on the basis of formula"
$ws.Range("E13").Value = "None"
$ws.Range("F13").Value = "account_number = 2000`nclient_number = 2000`nbalance = 20`ndate_created = (2015, 1, 1)`noverdraft_limit = -20`noverdraft_rate = 0.06"
$ws.Range("G13").Value = "service_charge set to base service charge "
$ws.Range("E14").Value = "None"
$ws.Range("F14").Value = "account_number = 2000`nclient_number = 2000`nbalance = 2000`ndate_created = (2015, 1, 1)`noverdraft_limit =- 20`noverdraft_rate = 0.06"
$ws.Range("G14").Value = "returned formatted string"

$ws.Rows.Item(7).RowHeight = 31.2
$ws.Rows.Item(8).RowHeight = 31.2
$ws.Rows.Item(9).RowHeight = 31.2
$ws.Rows.Item(10).RowHeight = 31.2
$ws.Rows.Item(11).RowHeight = 31.2
$ws.Rows.Item(12).RowHeight = 31.2
$ws.Rows.Item(13).RowHeight = 31.2

$ws.Range("C2:G2").Select()
